$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "JSU(-1.1142012660484941, 1.3176538011826278, 0.06027167841068085, 3.3895117290585137)"
$ws.Range("C2").Value = "JSU(-1.3233411354060733, 1.0184286301860583, 2.493018866458735, 3.6158934158163936)"
$ws.Range("D2").Value = "NIG(0.8723006976325841, 0.6010894403174518, 1.1616923223606455, 2.9423732139194296)"
$ws.Range("E2").Value = "JSU(-0.8811274555132748, 1.1522770564134794, 4.051790353589574, 5.250792110200804)"
